# Insert a new header row at the top of the 'CISPR 25' sheet, pushing
# all existing data down by one row, and populate the header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CISPR 25")

# Insert a new blank row above row 1; this shifts all data (and the
# named range CISPR_25, which Excel keeps in sync automatically) down.
$ws.Rows.Item(1).Insert()

# Fill in the new header row
$ws.Range("A1").Value = "F, Mhz"
$ws.Range("B1").Value = "Eeq max ref"

# Header row formatting: wrap text and taller row height
$ws.Range("A1:B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 30

# Update selection to match target state
$ws.Range("B2").Select()

# Ensure the CISPR_25 named range reflects the shifted data range
# (some environments do not auto-adjust defined names on row insert).
$wb.Names.Item("CISPR_25").RefersTo = "='CISPR 25'!`$A`$2:`$B`$482"

$wb.Save()
